$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Update header timestamp text
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 12:05"

# 2) Re-rank Oman above Nigeria and Marruecos (rows 58-60), keeping countries ordered by
#    total cases descending. Oman's refreshed numbers now exceed Nigeria/Marruecos, and the
#    old Nigeria/Marruecos rows shift down one position each.
$ws.Range("A58").Value = "Oman"
$ws.Range("B58").Value = 7770
$ws.Range("C58").Value = 513
$ws.Range("D58").Value = 1933
$ws.Range("E58").Value = 5801
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 36

$ws.Range("A59").Value = "Nigeria"
$ws.Range("B59").Value = 7526
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 2174
$ws.Range("E59").Value = 5131
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 221

$ws.Range("A60").Value = "Marruecos"
$ws.Range("B60").Value = 7406
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 4638
$ws.Range("E60").Value = 2570
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 198

# 3) Refresh case counts for Bielorrusia (row 26)
$ws.Range("B26").Value = 36198
$ws.Range("C26").Value = 954
$ws.Range("D26").Value = 14155
$ws.Range("E26").Value = 21844
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 199

# 4) Refresh case counts for Rumania (row 40)
$ws.Range("B40").Value = 18070
$ws.Range("C40").Value = 213
$ws.Range("D40").Value = 11399
$ws.Range("E40").Value = 5492

# 5) Refresh case counts for Albania (row 111)
$ws.Range("B111").Value = 998
$ws.Range("C111").Value = 9
$ws.Range("D111").Value = 789
$ws.Range("E111").Value = 177
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 32
